$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, row 2 ("dbExcel" column) changes from the Neo4j data filename
# to the new Manifest.csv filename.
$ws.Range("C2").Value = "TC01_Bento_E2E_Select-All-Add-To-Cart_Manifest.csv"

# View change: drop the frozen/scrolled "topLeftCell" and zoom the sheet to 70%.
$excel.ActiveWindow.Zoom = 70
